# Apply updated cryptocurrency market data (price + 1h volume change)
# to Sheet1, matching the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell reference + its new literal text value.
# "AsText" marks values that look numeric (e.g. "1.003", "216.06")
# so they must be forced to text -- otherwise Excel would silently
# reinterpret them as numbers and normalize/round them, corrupting
# the source-formatted strings scraped from the site.
$updates = @(
    @{ Ref = "D2"; Value = "25.846.90"; AsText = $false },
    @{ Ref = "E2"; Value = "  -0.34%  "; AsText = $false },
    @{ Ref = "D3"; Value = "1.636.79"; AsText = $false },
    @{ Ref = "E3"; Value = "  -0.12%  "; AsText = $false },
    @{ Ref = "E4"; Value = "  +0.05%  "; AsText = $false },
    @{ Ref = "D5"; Value = "216.06"; AsText = $true },
    @{ Ref = "E5"; Value = "  +0.63%  "; AsText = $false },
    @{ Ref = "D6"; Value = "0.5071"; AsText = $true },
    @{ Ref = "E6"; Value = "  -0.07%  "; AsText = $false },
    @{ Ref = "D7"; Value = "1.003"; AsText = $true },
    @{ Ref = "E7"; Value = "  +0.07%  "; AsText = $false },
    @{ Ref = "D8"; Value = "0.2579"; AsText = $true },
    @{ Ref = "E8"; Value = "  +0.13%  "; AsText = $false },
    @{ Ref = "D9"; Value = "0.06441"; AsText = $true },
    @{ Ref = "E9"; Value = "  +1.41%  "; AsText = $false },
    @{ Ref = "D10"; Value = "19.57"; AsText = $true },
    @{ Ref = "E10"; Value = "  -1.22%  "; AsText = $false },
    @{ Ref = "D11"; Value = "0.07789"; AsText = $true },
    @{ Ref = "E11"; Value = "  +0.77%  "; AsText = $false },
    @{ Ref = "D12"; Value = "4.284"; AsText = $true },
    @{ Ref = "E12"; Value = "  -0.17%  "; AsText = $false },
    @{ Ref = "B13"; Value = "WrappedEther"; AsText = $false },
    @{ Ref = "C13"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; AsText = $false },
    @{ Ref = "D13"; Value = "1.633.84"; AsText = $false },
    @{ Ref = "E13"; Value = "  -0.32%  "; AsText = $false },
    @{ Ref = "B14"; Value = "WrappedliquidstakedEther2.0"; AsText = $false },
    @{ Ref = "C14"; Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; AsText = $false },
    @{ Ref = "D14"; Value = "1.858.20"; AsText = $false },
    @{ Ref = "E14"; Value = "  -0.46%  "; AsText = $false },
    @{ Ref = "D15"; Value = "0.5616"; AsText = $true },
    @{ Ref = "E15"; Value = "  +2.83%  "; AsText = $false },
    @{ Ref = "B16"; Value = "ShibaInu"; AsText = $false },
    @{ Ref = "C16"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; AsText = $false },
    @{ Ref = "D16"; Value = "0.0₅7604"; AsText = $false },
    @{ Ref = "E16"; Value = "  -1.79%  "; AsText = $false },
    @{ Ref = "B17"; Value = "Litecoin"; AsText = $false },
    @{ Ref = "C17"; Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; AsText = $false },
    @{ Ref = "D17"; Value = "63.22"; AsText = $true },
    @{ Ref = "E17"; Value = "  -1.54%  "; AsText = $false },
    @{ Ref = "D18"; Value = "25.846.41"; AsText = $false },
    @{ Ref = "E18"; Value = "  -0.52%  "; AsText = $false },
    @{ Ref = "E19"; Value = "  +0.17%  "; AsText = $false },
    @{ Ref = "D20"; Value = "194.85"; AsText = $true },
    @{ Ref = "E20"; Value = "  -0.58%  "; AsText = $false },
    @{ Ref = "D21"; Value = "4.334"; AsText = $true },
    @{ Ref = "E21"; Value = "  -2.80%  "; AsText = $false },
    @{ Ref = "D22"; Value = "9.871"; AsText = $true },
    @{ Ref = "E22"; Value = "  -0.76%  "; AsText = $false },
    @{ Ref = "D23"; Value = "6.059"; AsText = $true },
    @{ Ref = "E23"; Value = "  -1.15%  "; AsText = $false },
    @{ Ref = "E24"; Value = "  +0.05%  "; AsText = $false },
    @{ Ref = "D25"; Value = "1.792"; AsText = $true },
    @{ Ref = "E25"; Value = "  -5.33%  "; AsText = $false },
    @{ Ref = "D26"; Value = "0.1281"; AsText = $true },
    @{ Ref = "E26"; Value = "  +1.27%  "; AsText = $false },
    @{ Ref = "D27"; Value = "140.58"; AsText = $true },
    @{ Ref = "D28"; Value = "6.766"; AsText = $true },
    @{ Ref = "E28"; Value = "  -1.32%  "; AsText = $false },
    @{ Ref = "D29"; Value = "15.48"; AsText = $true },
    @{ Ref = "E29"; Value = "  -1.02%  "; AsText = $false },
    @{ Ref = "D30"; Value = "1.241"; AsText = $true },
    @{ Ref = "E30"; Value = "  +0.28%  "; AsText = $false },
    @{ Ref = "D31"; Value = "0.04892"; AsText = $true },
    @{ Ref = "E31"; Value = "  -0.01%  "; AsText = $false },
    @{ Ref = "D32"; Value = "3.292"; AsText = $true },
    @{ Ref = "E32"; Value = "  +0.87%  "; AsText = $false },
    @{ Ref = "D33"; Value = "3.221"; AsText = $true },
    @{ Ref = "E33"; Value = "  +0.60%  "; AsText = $false },
    @{ Ref = "D34"; Value = "1.556"; AsText = $true },
    @{ Ref = "E34"; Value = "  +0.65%  "; AsText = $false },
    @{ Ref = "D35"; Value = "2.375"; AsText = $true },
    @{ Ref = "E35"; Value = "  +0.00%  "; AsText = $false },
    @{ Ref = "D36"; Value = "0.8993"; AsText = $true },
    @{ Ref = "E36"; Value = "  -1.86%  "; AsText = $false },
    @{ Ref = "D37"; Value = "2.569"; AsText = $true },
    @{ Ref = "D38"; Value = "1.129.55"; AsText = $false },
    @{ Ref = "E38"; Value = "  -0.11%  "; AsText = $false },
    @{ Ref = "D39"; Value = "0.5512"; AsText = $true },
    @{ Ref = "E39"; Value = "  -0.28%  "; AsText = $false },
    @{ Ref = "D40"; Value = "0.01560"; AsText = $true },
    @{ Ref = "E40"; Value = "  -0.55%  "; AsText = $false },
    @{ Ref = "D41"; Value = "0.9943"; AsText = $true },
    @{ Ref = "E41"; Value = "  -0.73%  "; AsText = $false },
    @{ Ref = "D42"; Value = "5.535"; AsText = $true },
    @{ Ref = "E42"; Value = "  -1.00%  "; AsText = $false },
    @{ Ref = "D43"; Value = "0.8001"; AsText = $true },
    @{ Ref = "E43"; Value = "  -0.32%  "; AsText = $false },
    @{ Ref = "D44"; Value = "97.32"; AsText = $true },
    @{ Ref = "E44"; Value = "  -1.34%  "; AsText = $false },
    @{ Ref = "D45"; Value = "1.783.67"; AsText = $false },
    @{ Ref = "E45"; Value = "  +0.33%  "; AsText = $false },
    @{ Ref = "E46"; Value = "  -7.60%  "; AsText = $false },
    @{ Ref = "B47"; Value = "Aave"; AsText = $false },
    @{ Ref = "C47"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; AsText = $false },
    @{ Ref = "D47"; Value = "55.54"; AsText = $true },
    @{ Ref = "E47"; Value = "  +0.60%  "; AsText = $false },
    @{ Ref = "B48"; Value = "Mantle"; AsText = $false },
    @{ Ref = "C48"; Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; AsText = $false },
    @{ Ref = "D48"; Value = "0.4440"; AsText = $true },
    @{ Ref = "E48"; Value = "  -1.81%  "; AsText = $false },
    @{ Ref = "D49"; Value = "0.05054"; AsText = $true },
    @{ Ref = "E49"; Value = "  -2.54%  "; AsText = $false },
    @{ Ref = "D50"; Value = "7.679"; AsText = $true },
    @{ Ref = "E50"; Value = "  +2.42%  "; AsText = $false },
    @{ Ref = "D51"; Value = "0.9997"; AsText = $true },
    @{ Ref = "E51"; Value = "  +0.05%  "; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.AsText) {
        # Pre-format as Text so the numeric-looking string is stored
        # verbatim, then drop the format override so the cell keeps
        # its original (default) style -- only the value changes.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
